$wb = $excel.ActiveWorkbook

# --- Rename "X"/"Y" parameter placeholder labels to "ParameterN" labels ---

$wsInt = $wb.Worksheets.Item("testM1Int")
$wsInt.Range("C2").Value = "getM1IntParameter0"
$wsInt.Range("C2").Select()

$wsFloat = $wb.Worksheets.Item("testM1Float")
$wsFloat.Range("C2").Value = "getM1FloatParameter0"
$wsFloat.Range("C2").Select()

$wsDouble = $wb.Worksheets.Item("testM1Double")
$wsDouble.Range("C2").Value = "getM1DoubleParameter0"
$wsDouble.Range("C2").Select()

$wsBoolean = $wb.Worksheets.Item("testM1Boolean")
$wsBoolean.Range("C2").Value = "getM1BooleanParameter0"
$wsBoolean.Range("C2").Select()

$wsString = $wb.Worksheets.Item("testM1String")
$wsString.Range("C2").Value = "getM1StringParameter0"
$wsString.Range("C2").Select()

$wsIntInt = $wb.Worksheets.Item("testM1IntInt")
$wsIntInt.Range("C2").Value = "getM1IntIntParameter0"
$wsIntInt.Range("D2").Value = "getM1IntIntParameter1"

# testM1IntInt becomes the active sheet/tab, with E14:E15 selected
$wsIntInt.Activate()
$wsIntInt.Range("E14:E15").Select()
